$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement dataset for rows 2-31 (x, y, z), reflecting:
#  - 3 new rows inserted at the top (new rows 2-4)
#  - original 20 rows shifted down (now rows 5-24)
#  - 7 new rows appended at the bottom (new rows 25-31)
$rows = @(
    @(0.3769224858049733, -0.2748469442259254, 0.004896030350543701),
    @(-0.8251284312733457, 0.638507986693405, 1.034206660098391),
    @(-0.384797418221631, 2.191421729508453, 1.820474159769618),
    @(-0.8565583720347749, 1.682360221204802, 1.74579396466501),
    @(-0.239309732830677, 0.5754852927407828, 1.096922339830855),
    @(0.1751386040748463, -0.2145548291081436, 0.8761335207905833),
    @(-1.060679824819472, 0.4943530759957171, 1.891958522744575),
    @(-0.450785677701105, -0.1085982450231191, 0.1413266023694169),
    @(0.3425711914571603, -2.118668488800253, -2.763528636449275),
    @(1.661536096485417, -2.177148738282218, -3.591206591202183),
    @(1.500267133962643, -2.098495937851318, -5.773556591641946),
    @(-2.984340389743152, -1.351582744756633, 1.246524475547487),
    @(-5.330522572629836, 1.363846176055872, 4.617774886335326),
    @(-4.297088194101645, 7.46291283436744, 1.662266135736299),
    @(-0.1566145409663497, 4.606385328884262, -0.8842834989056194),
    @(5.259676156606183, -3.144118193455617, -3.28624562092742),
    @(4.500417068015016, -0.31324529439601, -3.402522266171384),
    @(-1.445235498086886, 2.503525588710229, 0.09940216218542508),
    @(-8.082223704808689, 0.7967113996306803, 8.42470887134173),
    @(-1.0955425516487, -6.407353537588252, 3.211009324377832),
    @(3.420750588829331, -8.40992964615467, 0.8056608922617063),
    @(5.383590975182555, -4.346148073412959, -0.09347351565632778),
    @(1.778205815361046, -1.402306104331036, -0.230236690637843),
    @(-3.16472590750482, 0.5638766811924731, 1.501450692722371),
    @(-2.917193829232435, 0.08571250834320598, -0.5273152324309833),
    @(-0.5172856062260307, 4.862342070283832, -2.60658347866941),
    @(1.666729755276673, 0.2906178595197595, -1.515296798606113),
    @(2.873153332539526, -9.141991994266446, -1.962461360677373),
    @(1.324421677787197, -4.62868054152566, -0.05621947888198298),
    @(-3.513774491813999, 3.990525447645534, 6.551176387670226)
)

$rowCount = $rows.Count
$colCount = 3

$data = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$lastRow = 1 + $rowCount
$targetRange = $ws.Range("A2:C$lastRow")
$targetRange.Value2 = $data
